# Auto commit at 2026-01-16  7:54:22.97
# Update the Metrics sheet's raw input values (B2:B13). The "today" sheet
# pulls these through formulas (=Metrics!B2 etc.) and its own derived
# E/F columns, so those will recalculate automatically.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 284801.06
$metrics.Range("B3").Value  = 210035.44000000003
$metrics.Range("B4").Value  = 73714.640000000014
$metrics.Range("B5").Value  = 11642
$metrics.Range("B6").Value  = 5920671.7899999982
$metrics.Range("B7").Value  = 4980753.07
$metrics.Range("B8").Value  = 1737806.46
$metrics.Range("B9").Value  = 231919
$metrics.Range("B10").Value = 34386052.779999994
$metrics.Range("B11").Value = 32256028.23
$metrics.Range("B12").Value = 12019528.5
$metrics.Range("B13").Value = 1329549

# Restore the selection recorded on the Metrics sheet view.
$metrics.Range("D19").Select()

# The "today" sheet is the tab that was active/selected when the workbook
# was saved, with its own updated cell selection.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E7").Select()
